$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 230; this shifts the existing rows
# 230..333 down to 231..334 (matching the dimension growing from
# A1:R333 to A1:R334 seen in the diff).
$ws.Rows.Item(230).Insert()

# Populate the newly inserted row 230 with the new record.
$ws.Range("A230").Value = 6
$ws.Range("B230").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C230").Value = "Metropolitana"
$ws.Range("D230").Value = 45146
$ws.Range("E230").Value = 13
$ws.Range("F230").Value = 100112022
$ws.Range("G230").Value = "Arveja Verde"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 270
$ws.Range("K230").Value = 23000
$ws.Range("L230").Value = 25000
$ws.Range("M230").Value = 23889
$ws.Range("N230").Value = "$/malla 25 kilos"
$ws.Range("O230").Value = "Provincia de Limarí"
$ws.Range("P230").Value = 956
$ws.Range("Q230").Value = 25
$ws.Range("R230").Value = "Hortaliza"
